# Update the multiplication problems/answers in the table to match the
# regenerated output at commit c8c62b6.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "679×3=2037"; New = "155×7=1085" },
    @{ Old = "338×9=3042"; New = "760×5=3800" },
    @{ Old = "309×2=618";  New = "563×6=3378" },
    @{ Old = "536×6=3216"; New = "315×9=2835" },
    @{ Old = "199×6=1194"; New = "437×2=874"  },
    @{ Old = "580×9=5220"; New = "855×4=3420" },
    @{ Old = "725×9=6525"; New = "261×7=1827" },
    @{ Old = "430×6=2580"; New = "127×7=889"  },
    @{ Old = "460×7=3220"; New = "426×6=2556" },
    @{ Old = "767×9=6903"; New = "395×9=3555" },
    @{ Old = "350×2=700";  New = "954×3=2862" },
    @{ Old = "724×6=4344"; New = "823×4=3292" },
    @{ Old = "894×3=2682"; New = "212×7=1484" },
    @{ Old = "755×4=3020"; New = "167×7=1169" },
    @{ Old = "659×7=4613"; New = "921×3=2763" },
    @{ Old = "226×8=1808"; New = "422×6=2532" },
    @{ Old = "424×4=1696"; New = "302×7=2114" },
    @{ Old = "294×6=1764"; New = "454×9=4086" },
    @{ Old = "765×8=6120"; New = "124×9=1116" },
    @{ Old = "758×7=5306"; New = "957×9=8613" },
    @{ Old = "919×5=4595"; New = "349×7=2443" },
    @{ Old = "167×4=668";  New = "523×6=3138" },
    @{ Old = "435×3=1305"; New = "653×2=1306" },
    @{ Old = "717×7=5019"; New = "972×3=2916" },
    @{ Old = "232×2=464";  New = "923×2=1846" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $r.New, 2)
}

$d.Save()
